$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data block (new rows 2-4),
# shifting the existing data rows down by 3 (old row 2 becomes row 5, etc.)
$ws.Rows("2:4").Insert()

# The inserted rows pick up formatting from the row above (header); clear it
# so the new data cells remain unstyled, like the rest of the data rows.
$ws.Range("A2:C4").ClearFormats()

# Populate the newly inserted rows with their data (x, y, z)
$ws.Cells.Item(2, 1).Value = 0.01617096064405302
$ws.Cells.Item(2, 2).Value = -0.001510194632121166
$ws.Cells.Item(2, 3).Value = 0.006719517832001004
$ws.Cells.Item(3, 1).Value = 0.009587190579622892
$ws.Cells.Item(3, 2).Value = -0.01432139695518542
$ws.Cells.Item(3, 3).Value = 0.07086037078665362
$ws.Cells.Item(4, 1).Value = 0.0104356142692267
$ws.Cells.Item(4, 2).Value = 0.1928298026323316
$ws.Cells.Item(4, 3).Value = 0.05640322466691333

# Append 7 additional new data rows at the end of the sheet (rows 25-31)
$ws.Cells.Item(25, 1).Value = -0.168751522898674
$ws.Cells.Item(25, 2).Value = 0.04469497253497444
$ws.Cells.Item(25, 3).Value = -0.09234245866537093
$ws.Cells.Item(26, 1).Value = -0.1727900256713232
$ws.Cells.Item(26, 2).Value = -0.1038810287912686
$ws.Cells.Item(26, 3).Value = 0.06023810141616398
$ws.Cells.Item(27, 1).Value = -0.09510832776625952
$ws.Cells.Item(27, 2).Value = -0.07622240483760839
$ws.Cells.Item(27, 3).Value = 0.01844473597076204
$ws.Cells.Item(28, 1).Value = -0.0361937656998634
$ws.Cells.Item(28, 2).Value = 0.0348193198442459
$ws.Cells.Item(28, 3).Value = -0.0697913542389869
$ws.Cells.Item(29, 1).Value = -0.03700825323661167
$ws.Cells.Item(29, 2).Value = 0.04744386838542084
$ws.Cells.Item(29, 3).Value = 0.007177666657500789
$ws.Cells.Item(30, 1).Value = -0.0347344755298561
$ws.Cells.Item(30, 2).Value = 0.06149377011590525
$ws.Cells.Item(30, 3).Value = 0.0545706277092296
$ws.Cells.Item(31, 1).Value = -0.01959859269360693
$ws.Cells.Item(31, 2).Value = 0.00315613796313606
$ws.Cells.Item(31, 3).Value = 0.01527163075904057
